$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A4").Value = "bocaioandoru12+1@gmail.com"
$ws.Range("B4").Value = "Doru2"
$ws.Range("C4").Value = "Management"

$ws.Hyperlinks.Add($ws.Range("A4"), "mailto:bocaioandoru12+1@gmail.com")
$ws.Range("A4").Style = "Hyperlink"

$ws.Range("J13").Select()
